# Automatic update of files.
# Re-applies the refreshed case-list data: every "Förändrad" (C) date moves
# forward a day (46072 -> 46073), and the case rows (A/B/G, plus F where a
# "Kommuner" owner tag is present) are refreshed to the newly fetched order.
# The per-row hyperlink formulas (S/T/V/W/X/Y) are regenerated from the
# (possibly updated) case id (column A) so they keep pointing at the right
# per-case asset URLs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A="A 2943-2023"; B=44945; G=11.3; F=$null },
    @{ Row=3; A="A 18319-2025"; B=45762; G=3.4; F=$null },
    @{ Row=4; A="A 1782-2024"; B=45307; G=2.7; F=$null },
    @{ Row=5; A="A 38013-2022"; B=44811; G=5.3; F=$null },
    @{ Row=6; A="A 389-2023"; B=44929; G=2.5; F=$null },
    @{ Row=7; A="A 29245-2021"; B=44361; G=0.4; F=$null },
    @{ Row=8; A="A 48974-2023"; B=45209; G=4.5; F=$null },
    @{ Row=9; A="A 4493-2024"; B=45327; G=1.8; F=$null },
    @{ Row=10; A="A 4822-2023"; B=44957; G=2.2; F=$null },
    @{ Row=11; A="A 1531-2022"; B=44573; G=1.6; F=$null },
    @{ Row=12; A="A 4486-2024"; B=45327; G=0.6; F=$null },
    @{ Row=13; A="A 32610-2024"; B=45513; G=0.5; F=$null },
    @{ Row=14; A="A 10710-2025"; B=45722; G=1.8; F="Kommuner" },
    @{ Row=15; A="A 635-2023"; B=44930; G=0.5; F=$null },
    @{ Row=16; A="A 48181-2024"; B=45589; G=0.7; F=$null },
    @{ Row=17; A="A 34400-2025"; B=45846.61351851852; G=1.3; F=$null },
    @{ Row=18; A="A 34401-2025"; B=45846.6140162037; G=2.8; F=$null },
    @{ Row=19; A="A 18328-2025"; B=45762; G=1.8; F=$null },
    @{ Row=20; A="A 5817-2025"; B=45694.74113425926; G=1.2; F=$null },
    @{ Row=21; A="A 24-2023"; B=44928; G=0.5; F=$null },
    @{ Row=22; A="A 21572-2023"; B=45063; G=1.7; F=$null },
    @{ Row=23; A="A 7731-2026"; B=46062.52008101852; G=5.9; F=$null },
    @{ Row=24; A="A 7727-2026"; B=46062.50420138889; G=1.9; F=$null },
    @{ Row=25; A="A 28260-2023"; B=45099; G=5; F=$null },
    @{ Row=26; A="A 4481-2024"; B=45327; G=1; F=$null },
    @{ Row=27; A="A 18332-2025"; B=45762; G=2.5; F=$null },
    @{ Row=28; A="A 4256-2025"; B=45685; G=2; F=$null },
    @{ Row=29; A="A 11517-2024"; B=45372; G=0.7; F=$null },
    @{ Row=30; A="A 18434-2023"; B=45042; G=0.7; F=$null },
    @{ Row=31; A="A 53131-2021"; B=44468; G=1.3; F=$null },
    @{ Row=32; A="A 4487-2024"; B=45327; G=1.9; F=$null },
    @{ Row=33; A="A 18327-2025"; B=45762; G=0.6; F=$null }
)

# Rows whose hyperlink formulas (S/T/V/W/X/Y) embed the case id and must be
# regenerated to match a newly-assigned column-A value.
$relink = @(4, 6)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = 46073
    $ws.Range("G$row").Value = $r.G

    if ($r.F) {
        $ws.Range("F$row").Value = $r.F
    } else {
        $ws.Range("F$row").ClearContents()
    }

    if ($relink -contains $row) {
        $id = $r.A
        $ws.Range("S$row").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1278/artfynd/$id artfynd.xlsx`", `"$id`")"
        $ws.Range("T$row").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1278/kartor/$id karta.png`", `"$id`")"
        $ws.Range("V$row").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1278/klagomål/$id FSC-klagomål.docx`", `"$id`")"
        $ws.Range("W$row").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1278/klagomålsmail/$id FSC-klagomål mail.docx`", `"$id`")"
        $ws.Range("X$row").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1278/tillsyn/$id tillsynsbegäran.docx`", `"$id`")"
        $ws.Range("Y$row").Formula = "=HYPERLINK(`"https://klasma.github.io/Logging_1278/tillsynsmail/$id tillsynsbegäran mail.docx`", `"$id`")"
    }
}
